# Foundation properties.xlsx - "Add files via upload" edit
# Applies the data corrections / additions described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings, written in the same order they first appear in
#     the final workbook so the rebuilt shared-string table lines up with
#     the target ordering (K31, K27, F19, F27, E30/F30, F31). ---

# Row 31: Karagözyan Armenian Orphanage / Marriott Hotel and Key Plaza
#   new Contractor cell (K31) + corrected "Current state" (F31)
$ws.Range("K31").Value = "Bertuğ Bey İnşaat"

# Row 27: Inci Passage and Cinema / Lotus Nişantaşı
#   Contractor name typo fix (missing space) + corrected "Current state"
$ws.Range("K27").Value = "Bertuğ Bey İnşaat, Tabanlıoğlu Mimarlık"

# Row 19: Surp Agop apartment / Surp Agop Hospital polyclinic
#   "Current state" correction
$ws.Range("F19").Value = "polyclinic, offices"

# Row 27 "Current state" correction
$ws.Range("F27").Value = "offices, shopping center, residences"

# Row 30: Karagözyan Armenian Orphanage (cemetery row)
#   "Former state" / "Current state" correction
$ws.Range("E30").Value = "orphanage"
$ws.Range("F30").Value = "orphanage"

# Row 31 "Current state" correction
$ws.Range("F31").Value = "hotel, offices"

# --- New column M "spacer" cells with a thin boxed (medium grey L/R
#     border) style, matching rows that already show a Contractor value
#     (row 17 and row 22). ---

$m17 = $ws.Range("M17")
$m17.Value = " "
$m17.Borders.Weight = -4138
$m17.Borders.Color = 13421772
$m17.Borders.Item(8).LineStyle = -4142
$m17.Borders.Item(9).LineStyle = -4142

$m22 = $ws.Range("M22")
$m22.Value = " "
$m17.Copy()
$m22.PasteSpecial(-4122)

# --- Restore the view: clear the frozen/scrolled top-left cell and move
#     the active selection to L13 (matches the saved sheetView state). ---
$ws.Range("L13").Select()

Write-Host "Edit complete"
